$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Updated Diebold-Mariano statistics (column C) and p-values (column D)
$ws.Range("C2").Value = -0.1824824748765293
$ws.Range("D2").Value = 0.8568752146584846

$ws.Range("C3").Value = 0.01625958139480197
$ws.Range("D3").Value = 0.9871738608514755

$ws.Range("C4").Value = 1.599706644954175
$ws.Range("D4").Value = 0.1239274539652886

$ws.Range("C5").Value = 0.3678450420613426
$ws.Range("D5").Value = 0.7165041445181943

$ws.Range("C6").Value = 0.1791814761549157
$ws.Range("D6").Value = 0.8594351779687872

$ws.Range("C7").Value = 2.136190996932164
$ws.Range("D7").Value = 0.04403872731010927

$ws.Range("C8").Value = 0.4970158262641189
$ws.Range("D8").Value = 0.624108413987611

$ws.Range("C9").Value = 1.407887613110299
$ws.Range("D9").Value = 0.1731385250753494

$ws.Range("C10").Value = 0.326731706413469
$ws.Range("D10").Value = 0.7469589764799367

$ws.Range("C11").Value = -1.719075968560216
$ws.Range("D11").Value = 0.09964182829453527
